# Changed coding of significant column:
# Column B holds "RNAseq" (0/1) and column C holds "significant" (0/1).
# Whenever RNAseq (column B) is 0, the "significant" value in column C
# should no longer carry a 0/1 coding - it should be blank instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Data starts on row 2 (row 1 is the header row).
for ($r = 2; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    if ($b -eq 0) {
        $ws.Cells.Item($r, 3).ClearContents()
    }
}
